# Updates the "cryptos" price/volume table (Sheet1) to the latest scraped
# values. Most cells are plain strings (percentages, dotted price strings),
# so a direct .Value assignment is fine. A handful of Price cells look like
# plain decimals (e.g. "22.13") and Excel's COM layer would auto-coerce
# those to numbers (losing the original text formatting / trailing zeros),
# so for those we force the cell to Text format first, assign, then clear
# the number-format override again so the cell's style matches the rest of
# the column (no lingering "@" format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.975.35'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.556.58'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('E6').Value = '  +0.96%  '
$ws.Range('E7').Value = '  -0.58%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '22.13'
$c.ClearFormats()
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('E9').Value = '  -0.26%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0597'
$c.ClearFormats()
$ws.Range('E10').Value = '  +1.57%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0857'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '1.777.13'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '1.543.06'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '61.89'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.955.89'
$ws.Range('E17').Value = '  +0.19%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '218.01'
$c.ClearFormats()
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +2.74%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.32'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('E21').Value = '  -0.55%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.10'
$c.ClearFormats()
$ws.Range('E22').Value = '  +2.03%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.20'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.94%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.93'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.36%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '153.39'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('E26').Value = '  -0.16%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '15.02'
$c.ClearFormats()
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('E31').Value = '  -0.27%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.23'
$c.ClearFormats()
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.423.81'
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.12'
$c.ClearFormats()
$ws.Range('E34').Value = '  +4.24%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.08'
$c.ClearFormats()
$ws.Range('E35').Value = '  +13.44%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.ClearFormats()
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('E37').Value = '  +0.38%  '
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E39').Value = '  +2.11%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.809'
$c.ClearFormats()
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  -0.54%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.66'
$c.ClearFormats()
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('E44').Value = '  +0.89%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '64.68'
$c.ClearFormats()
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '1.690.90'
$ws.Range('E47').Value = '  +0.30%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '87.30'
$c.ClearFormats()
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('D50').Value = '0.0₇0991'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('E51').Value = '  +0.65%  '
